$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# Updated harmonized ratios (base_year change from #28 fix) recomputed for G3:G5
$ws.Range("G3").Value = 22.93333333333333
$ws.Range("G4").Value = 17.06666666666667
$ws.Range("G5").Value = 17.6
